# Update annotations for Ruilin
#
# The existing row 11 (Ruilin / 3 / 无 / DIS / MET / d3fb2dcb... / SJ3dBGZ0Z_annotated.xlsx / "We evaluate...")
# is kept, but its politeness_score (column B) becomes a real number (3) instead of text "3".
# A brand new row 12 is appended that duplicates the rest of that original row's data (still with
# politeness_score stored as text "3"), but carries a new id/source_file/text for a different annotation
# (afe80f3f... / r1BRfhiab_annotated.xlsx / "Not too surprisingly...").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new row 12: a copy of the original row 11 data, but pointing at the new annotation ---
$ws.Range("A12").Value = "Ruilin"
# Keep politeness_score as literal text "3" (matches the original row 11 formatting) by using the
# leading apostrophe so Excel stores it as text rather than coercing it to a number.
$ws.Range("B12").Value = "'3"
$ws.Range("C12").Value = "无"
$ws.Range("D12").Value = "DIS"
$ws.Range("E12").Value = "MET"
$ws.Range("F12").Value = "afe80f3f-3501-40b4-a3d0-1ad1f86c76ec"
$ws.Range("G12").Value = "r1BRfhiab_annotated.xlsx"
$ws.Range("H12").Value = "Not too surprisingly, the standard multiclass losses do not have the desired property, however approaches that reduce multi-class to binary classification at training time do, namely unnormalized models with penalized log Z (self-normalization), the NCE approach, as well as (the natural in the proposed setting) binary classification loss."

# --- Fix up row 11: only the politeness_score type changes, from text "3" to numeric 3 ---
$ws.Range("B11").Value = 3
